# Applies the cryptos list refresh described in the commit
# "Updated cryptos list on Wed May 22 15:11:48 UTC 2024 with GitHub Actions".
# All D/E value cells are plain text in the source workbook (t="inlineStr"),
# so any replacement that could be parsed by Excel as a number/date is
# written with a temporary "@" (text) number format, then the cell style is
# restored to "Normal" so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $range = $ws.Range($cell)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "70.117.37"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "3.742.62"
$ws.Range("E3").Value = "  -0.79%  "
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  -0.11%  "
Set-TextValue "D5" "617.32"
$ws.Range("E5").Value = "  +0.12%  "
Set-TextValue "D6" "181.70"
$ws.Range("E6").Value = "  +2.44%  "
$ws.Range("D7").Value = "3.745.89"
$ws.Range("E7").Value = "  -0.65%  "
Set-TextValue "D8" "1.00"
$ws.Range("E8").Value = "  -0.08%  "
Set-TextValue "D9" "0.533"
$ws.Range("E9").Value = "  -3.21%  "
$ws.Range("E10").Value = "  -1.42%  "
Set-TextValue "D11" "6.36"
$ws.Range("E11").Value = "  -0.69%  "
Set-TextValue "D12" "0.480"
$ws.Range("E12").Value = "  -4.94%  "
Set-TextValue "D13" "40.01"
$ws.Range("E13").Value = "  -1.41%  "
$ws.Range("E14").Value = "  -1.19%  "
$ws.Range("D15").Value = "4.362.50"
$ws.Range("E15").Value = "  -1.09%  "
$ws.Range("D16").Value = "3.737.55"
$ws.Range("E16").Value = "  -1.15%  "
$ws.Range("D17").Value = "70.154.13"
$ws.Range("E17").Value = "  -0.41%  "
$ws.Range("E18").Value = "  -2.09%  "
Set-TextValue "D19" "7.55"
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D20" "16.42"
$ws.Range("E20").Value = "  -2.70%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D21" "503.73"
$ws.Range("E21").Value = "  -4.03%  "
Set-TextValue "D22" "9.23"
$ws.Range("E22").Value = "  -2.89%  "
Set-TextValue "D23" "0.719"
$ws.Range("E23").Value = "  -3.69%  "
$ws.Range("E24").Value = "  +2.75%  "
Set-TextValue "D25" "86.87"
$ws.Range("E25").Value = "  -1.45%  "
Set-TextValue "D26" "12.93"
$ws.Range("E26").Value = "  -4.30%  "
Set-TextValue "D27" "11.18"
$ws.Range("E27").Value = "  +1.62%  "
Set-TextValue "D28" "0.0000131"
$ws.Range("E28").Value = "  +8.30%  "
$ws.Range("E29").Value = "  +0.03%  "
Set-TextValue "D30" "2.46"
$ws.Range("E30").Value = "  -2.29%  "
Set-TextValue "D31" "2.91"
$ws.Range("E31").Value = "  +0.83%  "
Set-TextValue "D32" "7.87"
$ws.Range("E32").Value = "  -1.27%  "
Set-TextValue "D33" "30.40"
$ws.Range("E33").Value = "  -5.64%  "
Set-TextValue "D34" "0.114"
$ws.Range("E34").Value = "  -0.85%  "
Set-TextValue "D35" "1.00"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("E36").Value = "  +0.49%  "
Set-TextValue "D37" "6.07"
$ws.Range("E37").Value = "  -1.55%  "
Set-TextValue "D38" "0.345"
$ws.Range("E38").Value = "  +0.82%  "
Set-TextValue "D39" "0.138"
$ws.Range("E39").Value = "  +3.56%  "
Set-TextValue "D40" "3.15"
$ws.Range("E40").Value = "  +13.61%  "
Set-TextValue "D41" "2.07"
$ws.Range("E41").Value = "  -4.75%  "
Set-TextValue "D42" "49.86"
$ws.Range("E42").Value = "  -3.30%  "
Set-TextValue "D43" "427.02"
$ws.Range("E43").Value = "  -0.09%  "
Set-TextValue "D44" "44.33"
$ws.Range("E44").Value = "  -0.67%  "
Set-TextValue "D45" "8.56"
$ws.Range("E45").Value = "  -3.68%  "
$ws.Range("D46").Value = "2.963.63"
$ws.Range("E46").Value = "  -5.59%  "
$ws.Range("E47").Value = "  -2.05%  "
Set-TextValue "D48" "27.12"
$ws.Range("E48").Value = "  -2.70%  "
Set-TextValue "D50" "136.24"
$ws.Range("E50").Value = "  -2.16%  "
$ws.Range("E51").Value = "  -2.65%  "
